$d = $word.ActiveDocument

function Find-ParaByText($text) {
    # Locate the (first) paragraph containing $text, anywhere in the body.
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if (-not $ok) {
        return $null
    }
    return $rng.Paragraphs(1)
}

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2) | Out-Null
}

function Set-ParaSpacing($para) {
    $para.Format.SpaceBefore = 3
    $para.Format.SpaceAfter = 0
    $para.Format.LineSpacingRule = 0
}

function Set-ParaRunSize($para, $size) {
    # Apply Font.Size to the paragraph's run text only (exclude the
    # trailing paragraph mark) so the paragraph-mark rPr is left untouched.
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    if ($end -gt $start) {
        $r = $d.Range($start, $end)
        $r.Font.Size = $size
    }
}

# ---------------------------------------------------------------
# Skills section
# ---------------------------------------------------------------

# "Programming languages" bullet - text unchanged, only add run sz=20
$p = Find-ParaByText("Programming languages")
Set-ParaRunSize $p 10

# "Software tools" bullet - text unchanged, add pPr spacing + run sz=20
$p = Find-ParaByText("Software tools")
Set-ParaSpacing $p
Set-ParaRunSize $p 10

# "Cloud and AI Tools" -> "Cloud and Data Management"
Replace-Text "Cloud and AI Tools: AWS, Azure, Google Cloud Platform, TensorFlow, and PyTorch" `
             "Cloud and Data Management: AWS, Azure, GCP basics"
$p = Find-ParaByText("Cloud and Data Management")
Set-ParaSpacing $p
Set-ParaRunSize $p 10

# "Cybersecurity" -> "AI/ML"
Replace-Text "Cybersecurity: Google Cybersecurity Professional Certificate, regularly participate in CTF challenges" `
             "AI/ML: Familiar with basic AI/ML concepts and integration"
$p = Find-ParaByText("AI/ML: Familiar")
Set-ParaSpacing $p
Set-ParaRunSize $p 10

# Insert new "Algorithm and Data Structures" bullet right after the AI/ML bullet
$newPara = $p.Range.InsertParagraphAfter()
$pAiMl = Find-ParaByText("AI/ML: Familiar")
$pNew = $pAiMl.Next()
$pNew.Range.Text = "• Algorithm and Data Structures: Experienced with university-level study"
Set-ParaSpacing $pNew
Set-ParaRunSize $pNew 10

# "Professional Skills" text update
Replace-Text "Professional Skills: Adaptability, Communication, Detail-oriented, Leadership, and Time Management" `
             "Professional Skills: Adaptable, Excellent communication, Detail-oriented, Leadership, Time Management"
$p = Find-ParaByText("Professional Skills")
Set-ParaSpacing $p
Set-ParaRunSize $p 10

# ---------------------------------------------------------------
# Experience - Job 1: Undergraduate Research Assistant
# ---------------------------------------------------------------

Replace-Text "❖ Undergraduate Research Assistant" "❖ Undergraduate Research Assistant (Node, React, JS)"

$pTitle = Find-ParaByText("Undergraduate Research Assistant")
Set-ParaRunSize $pTitle 10

$pSub = Find-ParaByText("University of Calgary, Calgary, AB")
Set-ParaSpacing $pSub
Set-ParaRunSize $pSub 10

Replace-Text "Developed automated data analysis workflows using Node.js and React, accelerating data extraction by 30%" `
             "Developed an automated workflow using Node and React for extracting detailed data insights in a timely manner."
$pB1 = Find-ParaByText("Developed an automated workflow")
Set-ParaSpacing $pB1
Set-ParaRunSize $pB1 10

Replace-Text "Integrated Cloud-based solutions to enhance information retrieval processes" `
             "Collected and processed multi-modal data (videos, spoken recordings, biometric data) for research in information needs."
$pB2 = Find-ParaByText("Collected and processed multi-modal data")
Set-ParaSpacing $pB2
Set-ParaRunSize $pB2 10

Replace-Text "Analyzed dataset patterns to propose AI-driven data processing improvements" `
             "Adapted quickly to new tools and technologies to enhance research data analysis processes, showcasing adaptability."
$pB3 = Find-ParaByText("Adapted quickly to new tools")
Set-ParaSpacing $pB3
Set-ParaRunSize $pB3 10

# Remove the 4th (now-trailing) bullet of job 1 entirely
$pB4 = Find-ParaByText("Collaborated with cross-functional teams to refine tools with real-time feedback")
$pB4.Range.Delete()

# ---------------------------------------------------------------
# Experience - Job 2: Self-Checkout Machine Software Developer -> Executive Team Member
# ---------------------------------------------------------------

$pJob2Title = Find-ParaByText("Self-Checkout Machine Software Developer")
Set-ParaSpacing $pJob2Title

Replace-Text "❖ Self-Checkout Machine Software Developer" "❖ Executive Team Member"
Replace-Text "Sep 2023 – Dec 2023" "Dec 2021 - Apr 2022"

$pJob2Title = Find-ParaByText("Executive Team Member")
Set-ParaRunSize $pJob2Title 10

$pJob2Sub = Find-ParaByText("Academic Project, University of Calgary")
Set-ParaSpacing $pJob2Sub

Replace-Text "Academic Project, University of Calgary" "Model United Nations (MUN) at Dar Jana International School"
$pJob2Sub = Find-ParaByText("Model United Nations (MUN) at Dar Jana International School")
Set-ParaRunSize $pJob2Sub 10

Replace-Text "Led a team of 20 in designing and developing Java software for self-checkout systems" `
             "Organized and prepared event documents, ensuring seamless execution of MUN conferences."
$pJ2B1 = Find-ParaByText("Organized and prepared event documents")
Set-ParaSpacing $pJ2B1
Set-ParaRunSize $pJ2B1 10

Replace-Text "Analyzed user interaction to optimize design for maximum transaction throughput" `
             "Managed participant engagement and facilitated communication, enhancing collaborative problem-solving skills."
$pJ2B2 = Find-ParaByText("Managed participant engagement")
Set-ParaSpacing $pJ2B2
Set-ParaRunSize $pJ2B2 10

Replace-Text "Applied data-driven feedback in interface design to enhance user experience" `
             "Demonstrated leadership by acting as a spokesman, guiding event procedures effectively."
$pJ2B3 = Find-ParaByText("Demonstrated leadership by acting as a spokesman")
Set-ParaSpacing $pJ2B3
Set-ParaRunSize $pJ2B3 10

# Remove the 4th (now-trailing) bullet of job 2 entirely
$pJ2B4 = Find-ParaByText("Built solutions aligned with security protocols to protect transaction data")
$pJ2B4.Range.Delete()

# ---------------------------------------------------------------
# Remove the entire "Full-stack Financial Assistant - Hackathon Project" block
# (title+date, subtitle, and its 4 bullets = 6 paragraphs)
# ---------------------------------------------------------------

$pBlockStart = Find-ParaByText("Full-stack Financial Assistant - Hackathon Project")
$pBlockEnd = Find-ParaByText("Analyzed feedback to iterate and enhance user interaction in real-time")
$blockRange = $d.Range($pBlockStart.Range.Start, $pBlockEnd.Range.End)
$blockRange.Delete()

Write-Output "done"
